$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update the time_taken column (F) with refreshed query timestamps
$ws.Range("F2").Value = "2021-10-05 14:21:33.405865"
$ws.Range("F3").Value = "2021-10-05 14:21:33.405873"
$ws.Range("F4").Value = "2021-10-05 14:21:33.405876"
$ws.Range("F5").Value = "2021-10-05 14:21:33.405878"
$ws.Range("F6").Value = "2021-10-05 14:21:33.405881"
$ws.Range("F7").Value = "2021-10-05 14:21:33.405885"
$ws.Range("F8").Value = "2021-10-05 14:21:33.405888"
$ws.Range("F9").Value = "2021-10-05 14:21:33.405890"
$ws.Range("F10").Value = "2021-10-05 14:21:33.405893"
$ws.Range("F11").Value = "2021-10-05 14:21:33.405895"
$ws.Range("F12").Value = "2021-10-05 14:21:33.405898"
$ws.Range("F13").Value = "2021-10-05 14:21:33.405900"
$ws.Range("F14").Value = "2021-10-05 14:21:33.405903"
$ws.Range("F15").Value = "2021-10-05 14:21:33.405905"
$ws.Range("F16").Value = "2021-10-05 14:21:33.405908"

# Add a new "metadata" tab after the "data" tab describing the panel query
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Mitochondrial disorder with complex III deficiency"
$meta.Range("C2").Value = 536
# Force "1.3" to stay a text value (matches the source "data_version" string) instead
# of being auto-coerced to a number: compute it as a text formula, then freeze the
# result down to a literal value so no stray number-format style is introduced.
$meta.Range("D2").Formula = '="1.3"'
$meta.Range("D2").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Range("E2").Value = "2020-02-17T16:02:10.171546Z"
$meta.Range("F2").Value = "2021-10-05 14:21:33.402559"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/536/?format=json"

# Match the header/index-column styling used on the "data" sheet (bold, thin border,
# centered horizontally, top-aligned vertically) by copying the existing cell format
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$ws.Activate()
